$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 132 (shifts old rows 132:227 down to 133:228,
# bringing their existing formatting/styles along for the ride).
$ws.Rows(132).Insert()

# Populate the newly-inserted row 132 with the new weekly record.
$ws.Range("A132").Value = 4
$ws.Range("B132").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C132").Value = "Los Lagos"
$ws.Range("D132").Value = 44651
$ws.Range("E132").Value = 10
$ws.Range("F132").Value = 100112017
$ws.Range("G132").Value = "Apio"
$ws.Range("H132").Value = "Americana (o)"
$ws.Range("I132").Value = "Primera"
$ws.Range("J132").Value = 25
$ws.Range("K132").Value = 12000
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = 12000
$ws.Range("N132").Value = "$/docena de matas"
$ws.Range("O132").Value = "Región de Coquimbo"
$ws.Range("P132").Value = 2000
$ws.Range("Q132").Value = 6
$ws.Range("R132").Value = "Hortaliza"
